$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample rows (3-10); only the header row and one data row remain.
$ws.Range("A3:D10").EntireRow.Delete()

# New header cells for the added columns
$ws.Range("E1").Value = "الصورة الشخصية"
$ws.Range("F1").Value = "الجنس"

# Replace the row 2 sample data with the real record
$ws.Range("A2").Value = "عبد المجيد"
$ws.Range("B2").Value = "الشامي"

# Keep the national number as literal text "9999" (not auto-converted to a number)
$ws.Range("C2").Formula = "=""9999"""
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

$ws.Range("D2").Value = "حمص"
$ws.Range("F2").Value = "ذكر"

# E2 holds a hyperlink ("show picture") styled blue + underlined
$ws.Range("E2").Value = "عرض الصورة"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://example.com")

# Drop the auto-generated named "Hyperlink" cell style so only a plain
# blue/underlined font + direct cell format remain (no extra named style).
$wb.Styles.Item("Hyperlink").Delete()
